# Regenerate merged AHB files
# - Rename the "_old"/"_new" header suffixes to "_FV2410"/"_FV2504"
# - Freeze the header row (row 1)
# - Turn the used range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename header row labels -----------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    # columns A..J (1..10) -> "<name>_FV2410"
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2410"
    # columns L..U (12..21) -> "<name>_FV2504"
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2504"
}
# column K (11) stays "diff" - unchanged

# --- 2. Freeze panes at row 2 (i.e. freeze header row 1) ------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into a table --------------------------------
$lastRow = $ws.Cells(1, 1).End([Microsoft.Office.Interop.Excel.XlDirection]::xlDown).Row
$lastCol = $ws.Cells(1, 1).End([Microsoft.Office.Interop.Excel.XlDirection]::xlToRight).Column
$tableRange = $ws.Range($ws.Cells(1, 1), $ws.Cells($lastRow, $lastCol))

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

Write-Host "Edit applied."
